$wb = $excel.ActiveWorkbook

# ---- Sheet: Schedule ----
$ws = $wb.Worksheets.Item("Schedule")

$ws.Range("E2").Value = 446.9438317499999
$ws.Range("F2").Value = 29.55977723214285

$ws.Range("E3").Value = -92.90290425000001
$ws.Range("F3").Value = -3.072185987103175

$ws.Range("E4").Value = 455.6571044999999
$ws.Range("F4").Value = 30.13605188492063

# ---- Sheet: Detailed ----
$ws = $wb.Worksheets.Item("Detailed")

$ws.Range("B9").Value = 57.06003
$ws.Range("B10").Value = 69.30083

$ws.Range("B11").Value = 67.02242
$ws.Range("C11").Value = "historical"

$ws.Range("B12").Value = 68.77846
$ws.Range("C12").Value = "historical"

$ws.Range("B13").Value = 78.54079
$ws.Range("B14").Value = 81.06431000000001

$ws.Range("B17").Value = 8.662850000000001

$ws.Range("B20").Value = -1.1692
$ws.Range("B21").Value = -5.58973
$ws.Range("B22").Value = -6.73663
$ws.Range("B23").Value = -6.89055
$ws.Range("B24").Value = -8.31794
$ws.Range("B25").Value = -7.69384
$ws.Range("B26").Value = -7.42697
$ws.Range("B27").Value = -7.38217
$ws.Range("B28").Value = -8.11548
$ws.Range("B29").Value = -7.76432
$ws.Range("B30").Value = -7.27164
$ws.Range("B31").Value = -7.31256
$ws.Range("B32").Value = -7.90483
$ws.Range("B33").Value = -6.41907
$ws.Range("B34").Value = -4.87048

$ws.Range("B36").Value = -5.01
$ws.Range("B37").Value = 4.45518
$ws.Range("B38").Value = 9.64795
$ws.Range("B39").Value = 33.31972
$ws.Range("B40").Value = 57.3
$ws.Range("B41").Value = 58.27097
$ws.Range("B42").Value = 59.91371
$ws.Range("B43").Value = 57.88115

$ws.Range("B45").Value = 56.98

$ws.Range("B47").Value = 58.29872
$ws.Range("B48").Value = 61.10806
$ws.Range("B49").Value = 58.79891
